$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match vs Rajasthan Royals (original row 2, "Oct 30 2020") and the
# match vs Sunrisers Hyderabad (original row 4, "Oct 24 2020") are both
# removed. The remaining match (vs Chennai Super Kings, "Nov 1 2020",
# originally row 3) shifts up to become the new row 2.
#
# Delete original row 2 first (Rajasthan Royals match). Everything below
# shifts up one row, so the former row 4 (Sunrisers Hyderabad match) is
# now row 3 - delete that next.
$ws.Rows(2).Delete()
$ws.Rows(3).Delete()
